# Generate Report for Handoff
# Updates the handoff identifiers, file hashes, and timestamps across the
# Overview / zh-cn / de-de sheets, keeping cell values and hyperlink
# "display" text in sync.

$wb = $excel.ActiveWorkbook

$oldGuid = "944f25c7-87e9-4061-94be-b8e994817dfa"
$newGuid = "db84d5ba-7972-4ed5-acda-9776b360321c"

$oldHash = "958064ad6fecfc4cf50a90b27882151cea2897de"
$newHash = "cfbc8dff5497595296480f3457e5da0a07012424"

$newMd  = "$newGuid.md"
$newZh  = "$newGuid.$newHash.zh-cn.xlf"
$newDe  = "$newGuid.$newHash.de-de.xlf"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value2 = $newMd
$wsOverview.Range("D2").Value2 = "2016-41-12 02:41:37"

foreach ($hh in $wsOverview.Hyperlinks) {
    if ($hh.Range.Address() -eq '$A$2') {
        $hh.TextToDisplay = $newMd
    }
}

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value2 = $newMd
$wsZh.Range("D2").Value2 = $newZh
$wsZh.Range("E2").Value2 = "2016-03-12 02:41:34"

foreach ($hh in $wsZh.Hyperlinks) {
    $addr = $hh.Range.Address()
    if ($addr -eq '$A$2') {
        $hh.TextToDisplay = $newMd
    } elseif ($addr -eq '$D$2') {
        $hh.TextToDisplay = $newZh
    }
}

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value2 = $newMd
$wsDe.Range("D2").Value2 = $newDe
$wsDe.Range("E2").Value2 = "2016-03-12 02:41:37"

foreach ($hh in $wsDe.Hyperlinks) {
    $addr = $hh.Range.Address()
    if ($addr -eq '$A$2') {
        $hh.TextToDisplay = $newMd
    } elseif ($addr -eq '$D$2') {
        $hh.TextToDisplay = $newDe
    }
}
